# tasks_14_12 - edited page's order and link and little edited last pages
#
# Adds the 14.12 work-log entry to row 24 (date, time in/out, time total and
# activity description) and updates the sheet's active selection to F25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 24: new "14.12" tracking entry -----------------------------------
$ws.Range("A24").Value = 14.12

$ws.Range("B24").Value = 0.375
$ws.Range("B24").NumberFormat = "h:mm"

$ws.Range("C24").Value = 0.58333333333333337
$ws.Range("C24").NumberFormat = "h:mm"

$ws.Range("E24").Value = "5hr"
$ws.Range("F24").Value = "trimed layout more and finish last page's layout and test"

# --- View state: scroll down one row and select F25 -----------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1

$ws.Range("F25").Select() | Out-Null
